$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '51.486.19'
$ws.Range('D3').Value = '2.786.09'
$ws.Range('E3').Value = '  -0.14%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '353.28'
$ws.Range('E5').Value = '  -1.44%  '
$ws.Range('D6').Value = '108.47'
$ws.Range('E6').Value = '  -1.14%  '
$ws.Range('D7').Value = '0.553'
$ws.Range('E7').Value = '  -1.25%  '
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').Value = '0.620'
$ws.Range('E9').Value = '  +5.09%  '
$ws.Range('D10').Value = '39.64'
$ws.Range('E10').Value = '  -1.79%  '
$ws.Range('E11').Value = '  +0.83%  '
$ws.Range('D12').Value = '0.0834'
$ws.Range('E12').Value = '  -1.90%  '
$ws.Range('D13').Value = '19.89'
$ws.Range('E13').Value = '  +1.85%  '
$ws.Range('D14').Value = '7.70'
$ws.Range('E14').Value = '  +1.38%  '
$ws.Range('D15').Value = '3.224.32'
$ws.Range('E15').Value = '  -0.03%  '
$ws.Range('D16').Value = '2.820.05'
$ws.Range('E16').Value = '  +2.00%  '
$ws.Range('D17').Value = '0.938'
$ws.Range('E17').Value = '  -1.23%  '
$ws.Range('D18').Value = '51.466.66'
$ws.Range('E18').Value = '  -0.79%  '
$ws.Range('D19').Value = '7.66'
$ws.Range('E19').Value = '  +2.86%  '
$ws.Range('D20').Value = '3.15'
$ws.Range('E20').Value = '  +2.28%  '
$ws.Range('D21').Value = '13.41'
$ws.Range('E21').Value = '  +1.58%  '
$ws.Range('D22').Value = '0.0₃0969'
$ws.Range('E22').Value = '  -0.77%  '
$ws.Range('D23').Value = '70.23'
$ws.Range('E23').Value = '  +0.18%  '
$ws.Range('D24').Value = '266.71'
$ws.Range('E24').Value = '  -1.47%  '
$ws.Range('D25').Value = '2.74'
$ws.Range('E25').Value = '  -0.81%  '
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  +0.06%  '
$ws.Range('D27').Value = '25.94'
$ws.Range('E27').Value = '  -2.02%  '
$ws.Range('E28').Value = '  -1.10%  '
$ws.Range('D29').Value = '10.29'
$ws.Range('E29').Value = '  +0.25%  '
$ws.Range('D30').Value = '36.54'
$ws.Range('E30').Value = '  +4.66%  '
$ws.Range('D31').Value = '6.14'
$ws.Range('E31').Value = '  +5.03%  '
$ws.Range('D32').Value = '51.90'
$ws.Range('E32').Value = '  -0.34%  '
$ws.Range('B33').Value = 'RenderToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D33').Value = '5.70'
$ws.Range('E33').Value = '  +9.47%  '
$ws.Range('B34').Value = 'VeChain'
$ws.Range('C34').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D34').Value = '0.0438'
$ws.Range('E34').Value = '  -6.12%  '
$ws.Range('B35').Value = 'Toncoin'
$ws.Range('C35').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D35').Value = '1.97'
$ws.Range('E35').Value = '  -8.53%  '
$ws.Range('D36').Value = '0.0851'
$ws.Range('E36').Value = '  +0.30%  '
$ws.Range('E37').Value = '  +0.05%  '
$ws.Range('D38').Value = '18.82'
$ws.Range('E38').Value = '  +0.25%  '
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').Value = '1.98'
$ws.Range('E39').Value = '  +0.11%  '
$ws.Range('B40').Value = 'LidoDAOToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D40').Value = '3.11'
$ws.Range('E40').Value = '  -3.03%  '
$ws.Range('E41').Value = '  +0.25%  '
$ws.Range('E42').Value = '  -4.53%  '
$ws.Range('D43').Value = '120.08'
$ws.Range('E43').Value = '  +0.68%  '
$ws.Range('E44').Value = '  -2.46%  '
$ws.Range('D45').Value = '21.59'
$ws.Range('E45').Value = '  -1.16%  '
$ws.Range('D46').Value = '2.118.24'
$ws.Range('E46').Value = '  +1.92%  '
$ws.Range('D47').Value = '3.35'
$ws.Range('E47').Value = '  +1.79%  '
$ws.Range('E48').Value = '  +5.31%  '
$ws.Range('B49').Value = 'SEI'
$ws.Range('C49').Value = 'https://coinranking.com/coin/8nxCqs-uj+sei-sei'
$ws.Range('D49').Value = '0.902'
$ws.Range('E49').Value = '  -5.09%  '
$ws.Range('B50').Value = 'TrustWalletToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D50').Value = '1.35'
$ws.Range('E50').Value = '  +7.78%  '
$ws.Range('D51').Value = '5.37'
$ws.Range('E51').Value = '  -7.41%  '
